$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 237 (pushing the
# existing rows 237..339 down to 238..340, which moves the former last
# row, 339, into the new row 340). The freshly inserted row 237 then
# receives what used to live in row 236 (same market entry, unchanged),
# while row 236 itself is overwritten with the brand-new observation
# (new date + new volume/price figures).

$ws.Rows.Item(237).Insert()

# Row 237 (new) = former row 236 contents, unchanged.
$ws.Range("A237").Value = 5
$ws.Range("B237").Value = "Macroferia Regional de Talca"
$ws.Range("C237").Value = "Maule"
$ws.Range("D237").Value = 44301
$ws.Range("E237").Value = 7
$ws.Range("F237").Value = 100114013
$ws.Range("G237").Value = "Zanahoria"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 500
$ws.Range("K237").Value = 6000
$ws.Range("L237").Value = 6000
$ws.Range("M237").Value = 6000
$ws.Range("N237").Value = "`$/saco 20 kilos"
$ws.Range("O237").Value = "Región de Ñuble"
$ws.Range("P237").Value = 300
$ws.Range("Q237").Value = 20
$ws.Range("R237").Value = "Hortaliza"

# Row 236 now holds the new observation (date + volume/prices changed).
$ws.Range("D236").Value = 44704
$ws.Range("J236").Value = 600
$ws.Range("K236").Value = 5500
$ws.Range("L236").Value = 5500
$ws.Range("M236").Value = 5500
$ws.Range("P236").Value = 275
